$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalSearch")

$ws.Range("C11").Value = '//*[@name="cmbChartOfAccount_input"]'
$ws.Range("A11").Value = "chartOfAccountInput"
$ws.Range("B11").Value = "by_xpath"

$ws.Range("A12").Value = "chartOfAccountList"
$ws.Range("D12").Value = '//*[@id="cmbChartOfAccount_listbox"]'
$ws.Range("C12").Value = '//*[@aria-controls="cmbChartOfAccount_listbox"]'
$ws.Range("B12").Value = "by_xpath"

$ws.Range("A13").Value = "leaseTypes"
$ws.Range("C13").Value = '//*[@aria-owns="ddlLeaseTypes_listbox"]'
$ws.Range("D13").Value = '//*[@id="ddlLeaseTypes_listbox"]'
$ws.Range("B13").Value = "by_xpath"

$ws.Columns.Item(3).ColumnWidth = 47.41

$wsLease = $wb.Worksheets.Item("Lease")
$wsLease.Activate()
$wsLease.Range("D2").Select()

$ws.Activate()
$ws.Range("D13").Select()
$ws.Application.ActiveWindow.ScrollRow = 4

